$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 584
$ws.Range("I2").Value = 584
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 584
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -471

$ws.Range("H38").Value = 168.85715
$ws.Range("I38").Value = 168.85715
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 506.57145
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -134.57145
$ws.Range("N38").ClearContents()

$ws.Range("H43").Value = 4595.4
$ws.Range("I43").Value = 2324.5
$ws.Range("J43").Value = 6109.3335
$ws.Range("K43").Value = 2324.5
$ws.Range("L43").Value = 6109.3335
$ws.Range("M43").Value = -2255.5
$ws.Range("N43").Value = -6247.3335

$ws.Range("H92").Value = 826.35297
$ws.Range("I92").Value = 660.2857
$ws.Range("J92").Value = 1601.3334
$ws.Range("K92").Value = 660.2857
$ws.Range("L92").Value = 1601.3334
$ws.Range("M92").Value = 587.7143

$ws.Range("H97").Value = 3632.25
$ws.Range("I97").Value = 8420
$ws.Range("J97").Value = 2036.3334
$ws.Range("K97").Value = 25260
$ws.Range("L97").Value = 6109.0002
$ws.Range("M97").Value = -24764
$ws.Range("N97").Value = -7101.0002

$ws.Range("H106").Value = 24434.428
$ws.Range("I106").Value = 24434.428
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 24434.428
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -23803.428

$ws.Range("H112").Value = 2249.8096
$ws.Range("I112").Value = 913
$ws.Range("J112").Value = 2472.611
$ws.Range("K112").Value = 2739
$ws.Range("L112").Value = 7417.833
$ws.Range("M112").Value = -1631
$ws.Range("N112").Value = -9633.832999999999

$ws.Range("H131").Value = 7062.375
$ws.Range("I131").Value = 999.5
$ws.Range("J131").Value = 9083.333000000001
$ws.Range("K131").Value = 2998.5
$ws.Range("L131").Value = 27249.999
$ws.Range("M131").Value = 2041.5

$ws.Range("H138").Value = 4161.196
$ws.Range("I138").Value = 1997.5
$ws.Range("J138").Value = 4367.2617
$ws.Range("K138").Value = 5992.5
$ws.Range("L138").Value = 13101.7851
$ws.Range("M138").Value = -852.5
$ws.Range("N138").Value = -23381.7851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2003.3334
$ws.Range("I2").Value = 2003.3334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2003.3334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1890.3334

$ws.Range("H32").Value = 7911.3955
$ws.Range("I32").Value = 5662.921
$ws.Range("J32").Value = 24999.8
$ws.Range("K32").Value = 5662.921
$ws.Range("L32").Value = 24999.8
$ws.Range("M32").Value = -5375.921

$ws.Range("H88").Value = 600
$ws.Range("I88").Value = 600
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 600
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -194
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 600
$ws.Range("I91").Value = 600
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 600
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 804
$ws.Range("N91").ClearContents()

$ws.Range("H116").Value = 2003.3334
$ws.Range("I116").Value = 2003.3334
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2003.3334
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 290.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2003.3334
$ws.Range("I3").Value = 2003.3334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2003.3334
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1889.3334

$ws.Range("H86").Value = 1508.7142
$ws.Range("I86").Value = 1512.2
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1512.2
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -389.2
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 1508.7142
$ws.Range("I89").Value = 1512.2
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 7561
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -1945
$ws.Range("N89").Value = -18732

$ws.Range("H94").Value = 3179.3845
$ws.Range("I94").Value = 2330.5
$ws.Range("J94").Value = 6009
$ws.Range("K94").Value = 2330.5
$ws.Range("L94").Value = 6009
$ws.Range("M94").Value = -1879.5

$ws.Range("H105").Value = 3831.75
$ws.Range("I105").Value = 2855.4285
$ws.Range("J105").Value = 10666
$ws.Range("K105").Value = 2855.4285
$ws.Range("L105").Value = 10666
$ws.Range("M105").Value = -1108.4285
$ws.Range("N105").Value = -14160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 17108.777
$ws.Range("I99").Value = 15287.111
$ws.Range("J99").Value = 18930.445
$ws.Range("K99").Value = 15287.111
$ws.Range("L99").Value = 18930.445
$ws.Range("M99").Value = -13789.111

$ws.Range("H122").Value = 5289
$ws.Range("I122").Value = 5610.8887
$ws.Range("J122").Value = 4323.3335
$ws.Range("K122").Value = 16832.6661
$ws.Range("L122").Value = 12970.0005
$ws.Range("M122").Value = -14382.6661
$ws.Range("N122").Value = -17870.0005

$ws.Range("H126").Value = 17108.777
$ws.Range("I126").Value = 15287.111
$ws.Range("J126").Value = 18930.445
$ws.Range("K126").Value = 45861.333
$ws.Range("L126").Value = 56791.335
$ws.Range("M126").Value = -43391.333

$ws.Range("H132").Value = 2172.8262
$ws.Range("I132").Value = 1377.8667
$ws.Range("J132").Value = 3663.375
$ws.Range("K132").Value = 4133.6001
$ws.Range("L132").Value = 10990.125
$ws.Range("M132").Value = -1603.6001

$ws.Range("H141").Value = 14999.333
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 14999.333
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 14999.333
$ws.Range("N141").Value = -25359.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125.833336
$ws.Range("I2").Value = 101.111115
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 606.66669
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -493.66669

$ws.Range("H34").Value = 1031
$ws.Range("I34").Value = 964
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 2892
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -2808
$ws.Range("N34").Value = -4668

$ws.Range("H55").Value = 2899.1667
$ws.Range("I55").Value = 798.3333
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 2394.9999
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -2217.9999
$ws.Range("N55").Value = -15354

$ws.Range("H80").Value = 4395.6665
$ws.Range("I80").Value = 2799
$ws.Range("J80").Value = 5992.3335
$ws.Range("K80").Value = 8397
$ws.Range("L80").Value = 17977.0005
$ws.Range("M80").Value = -7461

$ws.Range("H83").Value = 4395.6665
$ws.Range("I83").Value = 2799
$ws.Range("J83").Value = 5992.3335
$ws.Range("K83").Value = 25191
$ws.Range("L83").Value = 53931.0015
$ws.Range("M83").Value = -20511

$ws.Range("H113").Value = 1627.9
$ws.Range("I113").Value = 2333
$ws.Range("J113").Value = 1325.7142
$ws.Range("K113").Value = 6999
$ws.Range("L113").Value = 3977.1426
$ws.Range("M113").Value = -4829
$ws.Range("N113").Value = -8317.142599999999

$ws.Range("H131").Value = 1911.6666
$ws.Range("I131").Value = 1795
$ws.Range("J131").Value = 1970
$ws.Range("K131").Value = 5385
$ws.Range("L131").Value = 5910
$ws.Range("M131").Value = -345
$ws.Range("N131").Value = -15990

$ws.Range("H140").Value = 3843.3333
$ws.Range("I140").Value = 3843.3333
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 11529.9999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -6349.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7306.143
$ws.Range("I80").Value = 6035.75
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 6035.75
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = -5037.75
$ws.Range("N80").Value = -10996

$ws.Range("H83").Value = 7306.143
$ws.Range("I83").Value = 6035.75
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 30178.75
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -25186.75
$ws.Range("N83").Value = -54984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 60000
$ws.Range("I50").Value = 60000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 60000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -59363

$ws.Range("H55").Value = 480.8
$ws.Range("I55").Value = 440.125
$ws.Range("J55").Value = 643.5
$ws.Range("K55").Value = 440.125
$ws.Range("L55").Value = 643.5
$ws.Range("M55").Value = -267.125

$ws.Range("H61").Value = 3561.3572
$ws.Range("I61").Value = 4309
$ws.Range("J61").Value = 2564.5
$ws.Range("K61").Value = 4309
$ws.Range("L61").Value = 2564.5
$ws.Range("M61").Value = -4107
$ws.Range("N61").Value = -2968.5

$ws.Range("H113").Value = 3561.3572
$ws.Range("I113").Value = 4309
$ws.Range("J113").Value = 2564.5
$ws.Range("K113").Value = 4309
$ws.Range("L113").Value = 2564.5
$ws.Range("M113").Value = -2139
$ws.Range("N113").Value = -6904.5

$ws.Range("H122").Value = 8230.4
$ws.Range("I122").Value = 8230.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24691.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -22241.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 1000000
$ws.Range("I34").Value = 1000000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1000000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -999797

$ws.Range("H38").Value = 30000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 30000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30946
$ws.Range("M38").ClearContents()

$ws.Range("H81").Value = 19999.5
$ws.Range("I81").Value = 19999
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 39998
$ws.Range("L81").Value = 40000
$ws.Range("M81").Value = -38937

$ws.Range("H84").Value = 19999.5
$ws.Range("I84").Value = 19999
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 199990
$ws.Range("L84").Value = 200000
$ws.Range("M84").Value = -194686

$ws.Range("H100").Value = 1688.6154
$ws.Range("I100").Value = 1617
$ws.Range("J100").Value = 1849.75
$ws.Range("K100").Value = 3234
$ws.Range("L100").Value = 3699.5
$ws.Range("M100").Value = -2693

$ws.Range("H132").Value = 1870.3334
$ws.Range("I132").Value = 1668.2727
$ws.Range("J132").Value = 2426
$ws.Range("K132").Value = 5004.8181
$ws.Range("L132").Value = 7278
$ws.Range("M132").Value = -2474.8181
